$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The task ("RQ2") got a feedback pass - relabel the last existing entry
# (row 139) to reflect that, and log the new "Feedback einarbeiten" work
# on the following rows.
$ws.Range("C139").Value = "RQ2, Feedback einarbeiten"

$newRows = @(
    @{ Row = 140; Date = 45497; Hours = 8 },
    @{ Row = 141; Date = 45499; Hours = 4 },
    @{ Row = 142; Date = 45501; Hours = 1 },
    @{ Row = 143; Date = 45503; Hours = 6 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Date
    $ws.Cells.Item($r, 2).Value = $entry.Hours
    $ws.Cells.Item($r, 3).Value = "Feedback einarbeiten"

    # Match the formatting already used by the rows just above (row 139):
    # date format in column A, bold/black font style in columns B and C.
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Range("A139").NumberFormat
    $ws.Cells.Item($r, 2).Font.Color = $ws.Range("B139").Font.Color
    $ws.Cells.Item($r, 3).Font.Color = $ws.Range("C139").Font.Color
}

# Leave the view scrolled/selected near the newly added rows, same as the
# author's saved state after typing in the new entries.
$ws.Activate()
$ws.Range("A144").Select()
